$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-01-18"
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 21469
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 15259
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 36728
$ws.Range("L2").Value = 21471.1469
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 15260.5259
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 36731.6728

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-01-18"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 20621
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 14648
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 35269
$ws.Range("L3").Value = 20623.0621
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 14649.4648
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 35272.5269

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2025-01-18"
$ws.Range("A4").ClearFormats()
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 19828
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 13520
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 33348
$ws.Range("L4").Value = 19829.9828
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 13521.352
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 33351.3348

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2025-01-18"
$ws.Range("A5").ClearFormats()
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 19009
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 12691
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 31700
$ws.Range("L5").Value = 19010.9009
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 12692.2691
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 31703.17

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2025-01-18"
$ws.Range("A6").ClearFormats()
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 18488
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 12274
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 30762
$ws.Range("L6").Value = 18489.8488
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 12275.2274
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 30765.0762

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2025-01-18"
$ws.Range("A7").ClearFormats()
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 18753
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 12200
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 30953
$ws.Range("L7").Value = 18754.8753
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 12201.22
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 30956.0953

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2025-01-18"
$ws.Range("A8").ClearFormats()
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 18740
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 12590
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 31330
$ws.Range("L8").Value = 18741.874
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 12591.259
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 31333.133

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "2025-01-18"
$ws.Range("A9").ClearFormats()
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 19624
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 12971
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 32595
$ws.Range("L9").Value = 19625.9624
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 12972.2971
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 32598.2595

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "2025-01-18"
$ws.Range("A10").ClearFormats()
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 22373
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 13835
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 36208
$ws.Range("L10").Value = 22375.2373
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 13836.3835
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = 36211.6208

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "2025-01-18"
$ws.Range("A11").ClearFormats()
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 24582
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 14487
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 39069
$ws.Range("L11").Value = 24584.4582
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 14488.4487
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("T11").Value = 39072.9069

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2025-01-18"
$ws.Range("A12").ClearFormats()
$ws.Range("B12").Value = 11
$ws.Range("C12").Value = 25034
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 14850
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 39884
$ws.Range("L12").Value = 25036.5034
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 14851.485
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = 0
$ws.Range("S12").Value = 0
$ws.Range("T12").Value = 39887.9884

# Row 13
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "2025-01-18"
$ws.Range("A13").ClearFormats()
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 24749
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 14469
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 39218
$ws.Range("L13").Value = 24751.4749
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 14470.4469
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("T13").Value = 39221.9218

# Row 14
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "2025-01-18"
$ws.Range("A14").ClearFormats()
$ws.Range("B14").Value = 13
$ws.Range("C14").Value = 24350
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 14333
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 38683
$ws.Range("L14").Value = 24352.435
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 14334.4333
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("R14").Value = 0
$ws.Range("S14").Value = 0
$ws.Range("T14").Value = 38686.8683

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "2025-01-18"
$ws.Range("A15").ClearFormats()
$ws.Range("B15").Value = 14
$ws.Range("C15").Value = 25698
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 15180
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 40878
$ws.Range("L15").Value = 25700.5698
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 15181.518
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = 0
$ws.Range("S15").Value = 0
$ws.Range("T15").Value = 40882.0878

# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "2025-01-18"
$ws.Range("A16").ClearFormats()
$ws.Range("B16").Value = 15
$ws.Range("C16").Value = 26016
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 15492
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 41508
$ws.Range("L16").Value = 26018.6016
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 15493.5492
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = 0
$ws.Range("S16").Value = 0
$ws.Range("T16").Value = 41512.1508

# Row 17
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "2025-01-18"
$ws.Range("A17").ClearFormats()
$ws.Range("B17").Value = 16
$ws.Range("C17").Value = 25623
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 15143
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 40766
$ws.Range("L17").Value = 25625.5623
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 15144.5143
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = 0
$ws.Range("S17").Value = 0
$ws.Range("T17").Value = 40770.0766

# Row 18
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "2025-01-18"
$ws.Range("A18").ClearFormats()
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = 24213
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 14781
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 38994
$ws.Range("L18").Value = 24215.4213
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 14782.4781
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("R18").Value = 0
$ws.Range("S18").Value = 0
$ws.Range("T18").Value = 38997.8994

# Row 19
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "2025-01-18"
$ws.Range("A19").ClearFormats()
$ws.Range("B19").Value = 18
$ws.Range("C19").Value = 21283
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 14963
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 36246
$ws.Range("L19").Value = 21285.1283
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 14964.4963
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0
$ws.Range("R19").Value = 0
$ws.Range("S19").Value = 0
$ws.Range("T19").Value = 36249.6246
